$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (credentials, search filters) ---
$ws.Range("B20").Value = 'crypt:6dab1efafe436f246bfe392a2b64a16eb50f519d281360b3'
$ws.Range("B21").Value = 'crypt:adf6a09e5bcb826fd8eb2abadefa5770b5bace8a2ad9388a'
$ws.Range("B21").Font.Bold = $true
$ws.Range("B73").Value = 'Sale'
$ws.Range("B78").Value = '01-04-2023'
$ws.Range("B79").Value = '31-03-2024'

# --- Append new View-Register navigation rows (90-116) ---
$ws.Range("A90").Value = 'vr.No'
$ws.Range("B90").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[1]/div'
$ws.Range("A91").Value = 'vr.Inv.No'
$ws.Range("B91").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[2]/div'
$ws.Range("A92").Value = 'vr.Inv.Date'
$ws.Range("B92").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[3]/div'
$ws.Range("A93").Value = 'vr.Customer'
$ws.Range("B93").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[4]/div'
$ws.Range("A94").Value = 'vr.GSTIN'
$ws.Range("B94").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[5]/div'
$ws.Range("A95").Value = 'vr.TIN'
$ws.Range("B95").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[6]/div'
$ws.Range("A96").Value = 'vr.TaxFree'
$ws.Range("B96").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[7]/div'
$ws.Range("A97").Value = 'vr.Taxable'
$ws.Range("B97").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[8]/div'
$ws.Range("A98").Value = 'vr.GAmount'
$ws.Range("B98").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[9]/div'
$ws.Range("A99").Value = 'vr.table.row.xpath'
$ws.Range("B99").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/tbody/tr'
$ws.Range("A100").Value = 'vr.No.values'
$ws.Range("B100").Value = '/td[1]'
$ws.Range("A101").Value = 'vr.Inv.No.values'
$ws.Range("B101").Value = '/td[2]'
$ws.Range("A102").Value = 'vr.Inv.Date.values'
$ws.Range("B102").Value = '/td[3]'
$ws.Range("A103").Value = 'vr.Customer.values'
$ws.Range("B103").Value = '/td[4]'
$ws.Range("A104").Value = 'vr.GSTIN.values'
$ws.Range("B104").Value = '/td[5]'
$ws.Range("A105").Value = 'vr.TIN.values'
$ws.Range("B105").Value = '/td[6]'
$ws.Range("A106").Value = 'vr.TaxFree.values'
$ws.Range("B106").Value = '/td[7]'
$ws.Range("A107").Value = 'vr.Taxable.values'
$ws.Range("B107").Value = '/td[8]'
$ws.Range("A108").Value = 'vr.GAmount.values'
$ws.Range("B108").Value = '/td[9]'
$ws.Range("A109").Value = 'type.searchValue'
$ws.Range("B109").Value = '//main[@class=''mb-5'']/section/div[2]/section[2]/div/input'
$ws.Range("A110").Value = 'Search.Values.Register'
$ws.Range("B110").Value = '6/SL-24'
$ws.Range("A111").Value = 'present.search.values'
$ws.Range("B111").Value = '//tbody[@role=''rowgroup'']/tr/td[2]/a'
$ws.Range("A112").Value = 'empty.search'
$ws.Range("A113").Value = 'invoice.link'
$ws.Range("B113").Value = '//*[text()='' ${Search.Values.Register} '']'
$ws.Range("A114").Value = 'sales.Invoice.Number'
$ws.Range("B114").Value = '//span[text()='' 6/SL-24 '']'
$ws.Range("A115").Value = 'viewRegister.button'
$ws.Range("B115").Value = '//*[text()=''View Vouchers'']'
$ws.Range("A116").Value = 'voucher.presentvalues'
$ws.Range("B116").Value = '//*[text()=''Vouchers:'']'

# --- Move active selection to follow the newly appended content ---
$ws.Range("A117").Select()
